$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-11 down to 9-12
$ws.Rows.Item(8).Insert()

# Copy the date cell style (s="2") from the row above (now row 7) into new row 8,
# so the new date cell keeps the same date number format.
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new weekly record in row 8
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44729
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 412
$ws.Range("Q8").Value = 40
$ws.Range("R8").Value = "Hortaliza"
